$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header for the imported "tipo_cuenta" field
$ws.Range("G1").Value = "tipo_cuenta"

# Underlined placeholder cell used by the import/view logic
$ws.Range("G6").Font.Underline = $true

$ws.PageSetup.Orientation = 1
